# "falta unir el archivo excel del final con los docs"
#
# Update the SAP logon executable path (it moved from the old D:\...\ERPSAP
# install to the standard C:\...\SAP\FrontEnd install) and leave the
# workbook with the "Rutas" sheet active/selected at cell B3, as it was
# left by the author after editing that cell.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Rutas")
$ws2 = $wb.Worksheets.Item("parametrosInicio")

# Rutas!B2 : SAP logon path
$ws1.Range("B2").Value = "C:\Program Files (x86)\SAP\FrontEnd\SAPgui\saplogon.exe"

# Leave "Rutas" as the active sheet/tab with B3 selected.
$ws1.Activate()
$ws1.Range("B3").Select()
